# Weekly refresh of the "Fruta, Vega Modelo de Temuco - Caqui" data block.
# A new week's price row is inserted at row 6 (pushing the existing rows 6-12
# down to rows 7-13), and the oldest row (old row 13) is dropped, since rows
# 14+ are outside of this product's data window and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 <- new data
$ws.Range("D6").Value = 44685
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 1200

# Row 7 <- old row 6
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 17000
$ws.Range("Q7").Value = '$/bandeja 15 kilos granel'
$ws.Range("S7").Value = 1133
$ws.Range("T7").Value = 15

# Row 8 <- old row 7
$ws.Range("K8").Value = 'Fuyu'
$ws.Range("M8").Value = 120
$ws.Range("P8").Value = 10583
$ws.Range("S8").Value = 10583

# Row 9 <- old row 8
$ws.Range("D9").Value = 44305
$ws.Range("K9").Value = 'Mankaki'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 10500
$ws.Range("Q9").Value = '$/caja 10 kilos granel'
$ws.Range("S9").Value = 10500
$ws.Range("T9").Value = 1

# Row 10 <- old row 9
$ws.Range("D10").Value = 44312
$ws.Range("M10").Value = 135
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15481
$ws.Range("Q10").Value = '$/bandeja 15 kilos granel'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 1032

# Row 11 <- old row 10
$ws.Range("D11").Value = 44676
$ws.Range("K11").Value = 'Fuyu'
$ws.Range("M11").Value = 115
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = '$/bandeja 15 kilos'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1000

# Row 12 <- old row 11
$ws.Range("D12").Value = 44340
$ws.Range("K12").Value = 'Mankaki'
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 16471
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 1098

# Row 13 <- old row 12 (old row 13 data is discarded / falls out of the window)
$ws.Range("D13").Value = 44298
$ws.Range("M13").Value = 95
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("Q13").Value = '$/bandeja 15 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 667
$ws.Range("T13").Value = 15
